# Fixing this month balance problem, more advance copy function.
#
# 1. Remove the stray "h" cell (AF7) that was left over from a copy/paste,
#    which also drops the now-unused "h" shared string and re-collapses the
#    "Vacant Shifts" string's index (A7 keeps pointing at "Vacant Shifts").
# 2. Hide the helper/scratch columns AL:AQ (38-43) and the newly reclaimed
#    spacer column AR (44, width 0) used by the new "advance copy" helper
#    area, instead of deleting them outright.
# 3. Move the active selection to AH12 to match where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the leftover "h" value in AF7 - this also garbage-collects the now
# unused "h" shared string entry so "Vacant Shifts" (referenced by A7)
# shifts down to fill its slot.
$ws.Range("AF7").ClearContents() | Out-Null

# Hide columns AL:AQ (38-43), keeping their existing widths untouched.
$ws.Range("AL1:AQ1").EntireColumn.Hidden = $true

# Add/hide the new spacer column AR (44) with a zero display width.
# Excel's ColumnWidth setter always bakes in the ~0.8333 "padding" offset
# used for the default font, so subtract it back out to land on width=0.
$ws.Columns.Item(44).ColumnWidth = -0.8333333333333334
$ws.Columns.Item(44).Hidden = $true

# Move the saved selection to AH12.
$ws.Range("AH12").Select() | Out-Null
